# chore: adapt column header formatting to respective input file names (#7)
#
# The sheet compares an "old" AHB (Anwendungshandbuch) format version against
# a "new" one. The column headers used to carry generic "_old"/"_new" name
# suffixes; this adapts them to the concrete format-version identifiers of
# the two input files being diffed (FV2404 / FV2410). On top of that the
# data range is turned into a proper Excel Table and the header row is
# frozen so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells -------------------------------------------------
#   columns A:J   "<Name>_old" -> "<Name>_FV2404"
#   column  K     "diff"        (unchanged)
#   columns L:U   "<Name>_new" -> "<Name>_FV2410"

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($baseNames[$i] + "_FV2404")
}

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = ($baseNames[$i] + "_FV2410")
}

# --- 2. Wrap the data range in a native Excel Table -------------------------

$dataRange = $ws.UsedRange
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3. Freeze the header row so it stays visible while scrolling -----------

[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
